# Commit: Mon, May 04, 2020  9:06:27 AM
#
# Change the table style applied to the table on slide 6 (the
# "SOURCES OF FINANCE" table) from the deck's local custom table
# style ("Table_0") to the built-in PowerPoint table style
# {BD9FAA1E-CCEF-49C8-9EAC-B0B8D4C2FA43}.

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)

$table = $tableShape.Table
$table.ApplyStyle("{BD9FAA1E-CCEF-49C8-9EAC-B0B8D4C2FA43}")
